# edit.ps1
# Adds two new 96-well plate products (title / part number / image URL)
# to the products_URLimages worksheet, matching commit "added couple of wellplates".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column A to fit the longer titles being added ---
$ws.Columns.Item(1).ColumnWidth = 68.42578125

# --- new row 169 : 96-Well Plate, 2mL, Square Well, Round Bottom, Rim, Skirted, Raised Letters ---
# --- new row 170 : 96-Well Plate, 2mL, Square Well, V-Bottom, Raised Lettering, Sterile ---
# Fill column by column (A, then B, then C) so that shared-string insertion order
# mirrors how the workbook was actually authored.
$ws.Cells.Item(169, 1).Value = "96-Well Plate, 2mL, Square Well, Round Bottom, Rim, Skirted, Raised Letters"
$ws.Cells.Item(170, 1).Value = "96-Well Plate, 2mL, Square Well, V-Bottom, Raised Lettering, Sterile"

$ws.Cells.Item(169, 2).Value = "951652B"
$ws.Cells.Item(170, 2).Value = "951652C"
$ws.Cells.Item(169, 2).HorizontalAlignment = -4108
$ws.Cells.Item(170, 2).HorizontalAlignment = -4108

$ws.Cells.Item(169, 3).Value = "https://raw.githubusercontent.com/htslabs/images/main/951652B.jpg"
$ws.Cells.Item(170, 3).Value = "https://raw.githubusercontent.com/htslabs/images/main/951652C.jpg"

# --- hyperlink the two new image-url cells, matching the existing hyperlink style ---
$ws.Hyperlinks.Add($ws.Cells.Item(169, 3), "https://raw.githubusercontent.com/htslabs/images/main/951652B.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(170, 3), "https://raw.githubusercontent.com/htslabs/images/main/951652C.jpg") | Out-Null
$ws.Cells.Item(169, 3).Style = $ws.Range("C168").Style
$ws.Cells.Item(170, 3).Style = $ws.Range("C168").Style

# --- restore view state: scroll down to the new rows and select the next empty cell ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 155
$win.ScrollColumn = 1
$ws.Range("C172").Select()
